# Attic.docx: replace the inline "Attic guidelines" picture with a plain
# hyperlink run pointing at the image's original URL (the picture itself
# is dropped, the paragraph keeps its FirstParagraph style).

$d = $word.ActiveDocument

$url = "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Flats-Condominiums/F08_Attic.jpg?h=100%25&w=100%25"

# The picture is the sole InlineShape in the document (2nd paragraph).
$shape = $d.InlineShapes.Item(1)
$rng = $shape.Range

# Drop the picture but keep hold of its (now empty) range so the
# hyperlink gets inserted exactly where the drawing used to be.
$shape.Delete()

# Turn that spot into a hyperlink whose visible text is the URL itself.
$link = $d.Hyperlinks.Add($rng, $url, $null, $null, $url)

Write-Output "Replaced Attic picture with hyperlink: $($link.Address)"
